$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    ,@(10, "5:30 AM", 282.885, 312.885, 14.155)
    ,@(44, "5:47 AM", 1301.446667, 1331.446667, 14.786667)
    ,@(71, "6:00 AM", 2110.065, 2140.065, 8.81)
    ,@(81, "6:05 AM", 2418.395, 2448.395, 4.59)
    ,@(88, "6:09 AM", 2629.6425, 2659.6425, 5.215)
    ,@(105, "6:17 AM", 3136.858, 3166.858, 6.486)
    ,@(121, "6:25 AM", 3614.45, 3644.45, 3.595)
    ,@(146, "6:38 AM", 4379.04, 4409.04, 3.49)
    ,@(164, "6:47 AM", 4896.54, 4926.54, 4.39)
    ,@(174, "6:52 AM", 5193.04, 5223.04, 71.97)
    ,@(181, "6:56 AM", 5428.85, 5458.85, 1.735)
    ,@(222, "7:16 AM", 6648.33, 6678.33, 2.5925)
    ,@(252, "7:31 AM", 7544.065, 7574.065, 4.605)
    ,@(279, "7:44 AM", 8351.875, 8381.875, 0.06)
    ,@(290, "7:50 AM", 8698.795, 8728.795, 1.96)
    ,@(297, "7:53 AM", 8884.325000000001, 8914.325000000001, 4.97)
    ,@(306, "7:58 AM", 9166.17, 9196.17, 10.385)
    ,@(354, "8:22 AM", 10601.315, 10631.315, 0.745)
    ,@(367, "8:28 AM", 10997.11, 11027.11, 2.7)
    ,@(390, "8:40 AM", 11680.89, 11710.89, 0.763333)
    ,@(425, "8:57 AM", 12725.265, 12755.265, 4.76)
    ,@(451, "9:10 AM", 13511.86, 13541.86, 5.205)
    ,@(513, "9:41 AM", 15382.01, 15412.01, 0.865)
    ,@(532, "9:51 AM", 15943.655, 15973.655, 2.245)
    ,@(574, "10:12 A", 17202.985, 17232.985, 8.369999999999999)
    ,@(1014, "1:52 PM", 30408.62, 30438.62, 2.363333)
    ,@(1032, "2:01 PM", 30942.866667, 30972.866667, 9.316667000000001)
    ,@(1040, "2:05 PM", 31190.402, 31220.402, 0.65)
    ,@(1056, "2:13 PM", 31675.325, 31705.325, 6.43)
    ,@(1065, "2:17 PM", 31927.45, 31957.45, 12.815)
    ,@(1073, "2:21 PM", 32172.445, 32202.445, 12.6)
    ,@(1086, "2:28 PM", 32575.77, 32605.77, 21.33)
    ,@(1137, "2:53 PM", 34090.31, 34120.31, 4.425)
    ,@(1162, "3:06 PM", 34843.135, 34873.135, 1.29)
    ,@(1195, "3:22 PM", 35827.2775, 35857.2775, 4.8075)
    ,@(1212, "3:31 PM", 36344.033333, 36374.033333, 13.163333)
    ,@(1283, "4:06 PM", 38479.19, 38509.19, 51.24)
    ,@(1300, "4:15 PM", 38978.615, 39008.615, 27.355)
    ,@(1340, "4:35 PM", 40183.8725, 40213.8725, 9.875)
    ,@(1379, "4:54 PM", 41352.525, 41382.525, 1.775)
    ,@(1403, "5:06 PM", 42077.29, 42107.29, 0.585)
    ,@(1409, "5:10 PM", 42267.98, 42297.98, 7.51)
    ,@(1418, "5:14 PM", 42521.595, 42551.595, 25.965)
    ,@(1437, "5:23 PM", 43090.0675, 43120.0675, 3.4025)
    ,@(1472, "5:41 PM", 44138.865, 44168.865, 48.69)
    ,@(1484, "5:47 PM", 44500.423333, 44530.423333, 2.916667)
    ,@(1504, "5:57 PM", 45102.56, 45132.56, 3.275)
    ,@(1515, "6:02 PM", 45433.783333, 45463.783333, 3.276667)
    ,@(1536, "6:13 PM", 46066.956, 46096.956, 2.518)
    ,@(1569, "6:30 PM", 47068.13, 47098.13, 2.06)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
